# Fruta / hortaliza, semanal
# Insert two new weekly price records for "Clementina" (Mandarina) at the top
# of the date-ordered block (rows 269-270), pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above row 269; everything currently at 269.. shifts
# down to 271.. (formatting/styles are copied down from the row above, same
# as interactive Excel row-insert behaviour).
$ws.Rows("269:270").Insert()

# New row 269
$ws.Range("A269").Value = 5
$ws.Range("B269").Value = "Macroferia Regional de Talca"
$ws.Range("C269").Value = "Maule"
$ws.Range("D269").Value = 44753
$ws.Range("E269").Value = 7
$ws.Range("F269").Value = "Fruta"
$ws.Range("G269").Value = 100102
$ws.Range("H269").Value = "Cítricos"
$ws.Range("I269").Value = 100102004
$ws.Range("J269").Value = "Mandarina"
$ws.Range("K269").Value = "Clementina"
$ws.Range("L269").Value = "Especial"
$ws.Range("M269").Value = 150
$ws.Range("N269").Value = 8000
$ws.Range("O269").Value = 8000
$ws.Range("P269").Value = 8000
$ws.Range("Q269").Value = "$/caja 18 kilos"
$ws.Range("R269").Value = "Provincia de Quillota"
$ws.Range("S269").Value = 444
$ws.Range("T269").Value = 18

# New row 270
$ws.Range("A270").Value = 5
$ws.Range("B270").Value = "Macroferia Regional de Talca"
$ws.Range("C270").Value = "Maule"
$ws.Range("D270").Value = 44753
$ws.Range("E270").Value = 7
$ws.Range("F270").Value = "Fruta"
$ws.Range("G270").Value = 100102
$ws.Range("H270").Value = "Cítricos"
$ws.Range("I270").Value = 100102004
$ws.Range("J270").Value = "Mandarina"
$ws.Range("K270").Value = "Clementina"
$ws.Range("L270").Value = "Primera"
$ws.Range("M270").Value = 150
$ws.Range("N270").Value = 700
$ws.Range("O270").Value = 700
$ws.Range("P270").Value = 700
$ws.Range("Q270").Value = "$/caja 18 kilos"
$ws.Range("R270").Value = "Provincia de Quillota"
$ws.Range("S270").Value = 39
$ws.Range("T270").Value = 18
